$wb = $excel.ActiveWorkbook

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2770.4583
$ws.Range("I132").Value = 2825.6956
$ws.Range("K132").Value = 8477.086800000001
$ws.Range("M132").Value = -5947.086800000001

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 416.77777
$ws.Range("I135").Value = 195.14285
$ws.Range("K135").Value = 1756.28565
$ws.Range("M135").Value = 778.71435

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3421.5881
$ws.Range("I138").Value = 2046.25
$ws.Range("K138").Value = 6138.75
$ws.Range("M138").Value = -998.75

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2526.5715
$ws.Range("I63").Value = 1937.2
$ws.Range("K63").Value = 1937.2
$ws.Range("M63").Value = -1251.2

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2526.5715
$ws.Range("I66").Value = 1937.2
$ws.Range("K66").Value = 9686
$ws.Range("M66").Value = -6254

# ARM row 68
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 55000
$ws.Range("J68").Value = 55000
$ws.Range("L68").Value = 55000
$ws.Range("N68").Value = -56622

# ARM row 71
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H71").Value = 55000
$ws.Range("J71").Value = 55000
$ws.Range("L71").Value = 165000
$ws.Range("N71").Value = -173112

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4544.25
$ws.Range("I74").Value = 3655.3635
$ws.Range("K74").Value = 3655.3635
$ws.Range("M74").Value = -2781.3635

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4544.25
$ws.Range("I77").Value = 3655.3635
$ws.Range("K77").Value = 18276.8175
$ws.Range("M77").Value = -13908.8175

# ARM row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 39166.668
$ws.Range("J80").Value = 39166.668
$ws.Range("L80").Value = 39166.668
$ws.Range("N80").Value = -41162.668

# ARM row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 39166.668
$ws.Range("J83").Value = 39166.668
$ws.Range("L83").Value = 117500.004
$ws.Range("N83").Value = -127484.004

# ARM row 130
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H130").Value = 65429.332
$ws.Range("J130").Value = 65429.332
$ws.Range("L130").Value = 65429.332
$ws.Range("N130").Value = -75469.33199999999

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3718.5
$ws.Range("I134").Value = 3988
$ws.Range("J134").Value = 2999.8333
$ws.Range("K134").Value = 11964
$ws.Range("L134").Value = 8999.499899999999
$ws.Range("M134").Value = -9429
$ws.Range("N134").Value = -14069.4999

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2284.1272
$ws.Range("I31").Value = 1846.5526
$ws.Range("J31").Value = 3262.2354
$ws.Range("K31").Value = 1846.5526
$ws.Range("L31").Value = 3262.2354
$ws.Range("M31").Value = -1551.5526
$ws.Range("N31").Value = -3852.2354

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2284.1272
$ws.Range("I34").Value = 1846.5526
$ws.Range("J34").Value = 3262.2354
$ws.Range("K34").Value = 1846.5526
$ws.Range("L34").Value = 3262.2354
$ws.Range("M34").Value = -1644.5526
$ws.Range("N34").Value = -3666.2354

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7510.4443
$ws.Range("I99").Value = 7513.5713
$ws.Range("K99").Value = 7513.5713
$ws.Range("M99").Value = -6015.5713

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7510.4443
$ws.Range("I126").Value = 7513.5713
$ws.Range("K126").Value = 22540.7139
$ws.Range("M126").Value = -20070.7139

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2073
$ws.Range("I132").Value = 1502.2
$ws.Range("K132").Value = 4506.6
$ws.Range("M132").Value = -1976.6

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2065.4285
$ws.Range("I4").Value = 2039.2142
$ws.Range("K4").Value = 6117.642599999999
$ws.Range("M4").Value = -6005.642599999999

# CUL row 6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 5304.1665
$ws.Range("I6").Value = 108.333336
$ws.Range("J6").Value = 10500
$ws.Range("K6").Value = 325.000008
$ws.Range("L6").Value = 31500
$ws.Range("M6").Value = -212.000008
$ws.Range("N6").Value = -31726

# CUL row 25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1469.5
$ws.Range("I25").Value = 293
$ws.Range("K25").Value = 879
$ws.Range("M25").Value = -710

# CUL row 29
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 310
$ws.Range("I29").Value = 280
$ws.Range("J29").Value = 325
$ws.Range("K29").Value = 840
$ws.Range("L29").Value = 975
$ws.Range("M29").Value = -563
$ws.Range("N29").Value = -1529

# CUL row 30
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 1469.5
$ws.Range("I30").Value = 293
$ws.Range("K30").Value = 879
$ws.Range("M30").Value = -777

# CUL row 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 934.6667
$ws.Range("J51").Value = 1005
$ws.Range("L51").Value = 3015
$ws.Range("N51").Value = -3935

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1392
$ws.Range("J68").Value = 1225
$ws.Range("L68").Value = 3675
$ws.Range("N68").Value = -5297

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1392
$ws.Range("J71").Value = 1225
$ws.Range("L71").Value = 11025
$ws.Range("N71").Value = -19137

# CUL row 86
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 436.8
$ws.Range("I86").Value = 436.8
$ws.Range("K86").Value = 1310.4
$ws.Range("M86").Value = -124.4000000000001

# CUL row 89
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 436.8
$ws.Range("I89").Value = 436.8
$ws.Range("K89").Value = 3931.2
$ws.Range("M89").Value = 1996.8

# CUL row 138
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 6177.2856
$ws.Range("J138").Value = 6950.3335
$ws.Range("L138").Value = 20851.0005
$ws.Range("N138").Value = -31131.0005

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1550
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1550
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1550
$ws.Range("N97").Value = -2542
$ws.Range("M97").ClearContents()

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4006
$ws.Range("I132").Value = 4006
$ws.Range("K132").Value = 12018
$ws.Range("M132").Value = -9488

# LTW row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

# WVR row 20
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 886.5
$ws.Range("I126").Value = 984.1
$ws.Range("J126").Value = 398.5
$ws.Range("K126").Value = 2952.3
$ws.Range("L126").Value = 1195.5
$ws.Range("M126").Value = -482.3000000000002
$ws.Range("N126").Value = -6135.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1669.1538
$ws.Range("I132").Value = 1669.1538
$ws.Range("K132").Value = 5007.4614
$ws.Range("M132").Value = -2477.4614

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9460.200000000001
$ws.Range("I136").Value = 9971.714
$ws.Range("K136").Value = 29915.142
$ws.Range("M136").Value = -27365.142
